$d = $word.ActiveDocument

# --- 1. Rename the 2011 case-study source line -----------------------------
$old1 = "Files\\2011 Case Study\\Primary Sources_Policy_Strategies\\2010_national_security_strategy"
$new1 = "Files\\2011 Case Study\\CS1_Primary Sources_Policy_Strategies\\2010 National Security Strategy"
$d.Content.Find.Execute($old1, $true, $false, $false, $false, $false, $true, 1, $false, $new1, 2) | Out-Null

# --- 2. Rename the 2015 case-study source line -----------------------------
$old2 = "Files\\2015 Case Study\\Primary Sources_Policy_Strategies\\2015 National Security Strategy CLEAN"
$new2 = "Files\\2015 Case Study\\CS2_Primary Sources_Policy_Strategies\\2015 National Security Strategy"
$d.Content.Find.Execute($old2, $true, $false, $false, $false, $false, $true, 1, $false, $new2, 2) | Out-Null

# --- 3. Append the new 2018/2017 case-study block at the end of the doc ----
$wNs = 'xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"'

function Add-HighlightPara {
    param([string]$text)
    $last = $d.Paragraphs.Last
    $last.Range.InsertParagraphAfter()
    $p = $d.Paragraphs.Last
    $xml = "<w:p $wNs><w:pPr><w:pStyle w:val=""TextBody""/><w:bidi w:val=""0""/><w:spacing w:before=""113"" w:after=""113""/><w:ind w:left=""113"" w:right=""113"" w:hanging=""0""/><w:jc w:val=""left""/><w:rPr><w:highlight w:val=""lightGray""/></w:rPr></w:pPr><w:r><w:rPr><w:highlight w:val=""lightGray""/></w:rPr><w:t>$text</w:t></w:r></w:p>"
    $p.Range.InsertXML($xml)
    $p = $d.Paragraphs.Last
    $p.SpaceBefore = 5.65
    $p.SpaceAfter = 5.65
    $p.LeftIndent = 5.65
    $p.RightIndent = 5.65
}

function Add-PlainPara {
    param([string]$innerXml)
    $last = $d.Paragraphs.Last
    $last.Range.InsertParagraphAfter()
    $p = $d.Paragraphs.Last
    $xml = "<w:p $wNs><w:pPr><w:pStyle w:val=""TextBody""/><w:bidi w:val=""0""/><w:spacing w:before=""0"" w:after=""0""/><w:jc w:val=""left""/><w:rPr/></w:pPr><w:r><w:rPr/>$innerXml</w:r></w:p>"
    $p.Range.InsertXML($xml)
    $p = $d.Paragraphs.Last
    $p.SpaceBefore = 0
    $p.SpaceAfter = 0
}

Add-HighlightPara "Files\\2018 Case Study\\CS3_Primary Sources_Policy_Strategies\\2017 National Security Strategy - § 4 references coded [ 0.17% Coverage]"
Add-HighlightPara "Reference 1 - 0.04% Coverage"
Add-PlainPara '<w:t>e United States faces an extraordinarily dangerous world, &#xFB01; lled with a wide range of threats that have intensified in recent years.</w:t>'
Add-HighlightPara "Reference 2 - 0.04% Coverage"
Add-PlainPara '<w:t xml:space="preserve">A Competitive World The United States will respond to the growing </w:t><w:br/><w:t>political, economic, and military competitions we face around the world.</w:t>'
Add-HighlightPara "Reference 3 - 0.03% Coverage"
Add-PlainPara '<w:t xml:space="preserve">A </w:t><w:br/><w:t>central continuity in history is the contest for power. The present time period is no different.</w:t>'
Add-HighlightPara "Reference 4 - 0.06% Coverage"
Add-PlainPara '<w:t>The contests over influence are timeless. They have existed in varying degrees and levels of intensity, for millennia. Geopolitics is the interplay of these contests across the globe.</w:t>'

Write-Output "done: paragraphs=$($d.Paragraphs.Count)"
